$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New bucket-label values for columns T (X_luokka) and U (Y_luokka), rows 2-48.
# Values are the same text as before but re-typed with a leading space for
# every bucket whose starting number has fewer than 5 digits (Excel padding).
$labels = @(
    @(" 3000-3999", " 0-999"),
    @(" 2000-2999", " 2000-2999"),
    @(" 2000-2999", " 2000-2999"),
    @(" 2000-2999", " 2000-2999"),
    @(" 1000-1999", " 2000-2999"),
    @(" 1000-1999", " 2000-2999"),
    @(" 8000-8999", " 2000-2999"),
    @(" 4000-4999", " 2000-2999"),
    @(" 1000-1999", " 1000-1999"),
    @(" 8000-8999", " 2000-2999"),
    @(" 3000-3999", " 0-999"),
    @("13000-13999", " 3000-3999"),
    @(" 2000-2999", " 2000-2999"),
    @(" 2000-2999", " 2000-2999"),
    @(" 4000-4999", " 2000-2999"),
    @(" 4000-4999", " 2000-2999"),
    @(" 2000-2999", " 2000-2999"),
    @(" 9000-9999", " 2000-2999"),
    @("10000-10999", " 2000-2999"),
    @("10000-10999", " 2000-2999"),
    @(" 9000-9999", " 2000-2999"),
    @(" 2000-2999", " 2000-2999"),
    @(" 9000-9999", " 2000-2999"),
    @(" 1000-1999", " 1000-1999"),
    @(" 1000-1999", " 1000-1999"),
    @(" 8000-8999", " 2000-2999"),
    @(" 2000-2999", " 2000-2999"),
    @(" 8000-8999", " 2000-2999"),
    @(" 8000-8999", " 2000-2999"),
    @(" 1000-1999", " 3000-3999"),
    @(" 2000-2999", " 2000-2999"),
    @(" 3000-3999", " 0-999"),
    @(" 3000-3999", " 0-999"),
    @(" 9000-9999", " 2000-2999"),
    @(" 3000-3999", " 2000-2999"),
    @(" 3000-3999", " 2000-2999"),
    @("10000-10999", " 2000-2999"),
    @(" 8000-8999", " 2000-2999"),
    @(" 8000-8999", " 2000-2999"),
    @(" 8000-8999", " 2000-2999"),
    @(" 9000-9999", " 2000-2999"),
    @("10000-10999", " 2000-2999"),
    @(" 8000-8999", " 0-999"),
    @(" 8000-8999", " 2000-2999"),
    @(" 9000-9999", " 2000-2999"),
    @(" 2000-2999", " 2000-2999"),
    @(" 9000-9999", " 2000-2999")
)

for ($i = 0; $i -lt $labels.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 20).Value = $labels[$i][0]
    $ws.Cells.Item($r, 21).Value = $labels[$i][1]
}

# Column R (X_KA_ero) got manually widened.
$ws.Columns.Item(18).ColumnWidth = 17.140625

# Turn the data range into an AutoFilter table and select it.
$ws.Range("A1:U48").AutoFilter() | Out-Null
$ws.Range("A2:U48").Select() | Out-Null
